# Update the public EPEX Spot price workbook with the latest day of data.
#
# Sheet "Prix Spot": add a new date column BM ("17-aug") with hourly prices.
# Sheet "Gaz" and "CO2": append a new row for date 2025-08-15.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Prix Spot" -> add column BM (17-aug)
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Column BL is the last existing column (16-aug); BM is the new one (col 65).
$lastCol = 64   # BL
$newCol  = 65   # BM

# Copy the header cell's formatting (bold, border, centered) from BL1 to BM1,
# then overwrite its value - this keeps the new header cell styled exactly
# like the other date headers.
$srcHeader = $wsPrix.Cells.Item(1, $lastCol)
$dstHeader = $wsPrix.Cells.Item(1, $newCol)
$srcHeader.Copy($dstHeader)
$dstHeader.Value = "17-aug"

# Hourly values for the new day (rows 2..25).
$bmValues = @{
    2  = 67.88
    3  = 60.56
    4  = 53.86
    5  = 41.73
    6  = 37.24
    7  = 27.57
    8  = 29.12
    9  = 17.57
    10 = 30.65
    11 = 23.68
    12 = 9.699999999999999
    13 = 0.65
    14 = 0.65
    15 = 0
    16 = 0
    17 = 0.65
    18 = 5.13
    19 = 16.72
    20 = 51.45
    21 = 94.03
    22 = 108.6
    23 = 108.05
    24 = 101.82
    25 = 95.84
}

foreach ($row in 2..25) {
    $wsPrix.Cells.Item($row, $newCol).Value = $bmValues[$row]
}

# ---------------------------------------------------------------------
# Sheet "Gaz" -> append row 62 for date 2025-08-15
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$gazDateCell = $wsGaz.Cells.Item(62, 1)
# Force text format before assigning so the "YYYY-MM-DD" string is not
# auto-converted into a date serial number by Excel's input parsing.
$gazDateCell.NumberFormat = "@"
$gazDateCell.Value = "2025-08-15"
# Restore the default (unstyled) cell style, matching the other date cells
# in the column which carry no explicit formatting.
$gazDateCell.Style = "Normal"
$wsGaz.Cells.Item(62, 2).Value = 29.825

# ---------------------------------------------------------------------
# Sheet "CO2" -> append row 62 for date 2025-08-15
# ---------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")
$co2DateCell = $wsCO2.Cells.Item(62, 1)
$co2DateCell.NumberFormat = "@"
$co2DateCell.Value = "2025-08-15"
$co2DateCell.Style = "Normal"
$wsCO2.Cells.Item(62, 2).Value = 69.95
